$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '28.537.47'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'" + '  +0.56%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'" + '1.873.53'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'" + '  -0.06%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'" + '  -0.28%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'" + '315.70'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'" + '  +0.05%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'" + '  -0.44%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'" + '0.5075'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'" + '  -1.08%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'" + '0.3893'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'" + '  -0.66%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'" + '0.08346'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'" + '  +0.69%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("B10").Value = "'" + 'Polygon'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").Value = "'" + 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").Value = "'" + '1.102'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'" + '  -1.59%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("B11").Value = "'" + 'OKB'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").Value = "'" + 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").Value = "'" + '41.78'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'" + '  -0.27%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'" + '6.217'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'" + '  -0.68%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'" + '1.870.49'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'" + '  -1.04%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'" + '20.44'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'" + '  +0.44%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'" + '7.235'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'" + '  -0.11%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'" + '1.009'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'" + '  -0.23%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'" + '0.00001103'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'" + '  -0.05%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'" + '91.12'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'" + '  -0.21%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'" + '0.06705'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'" + '  -0.20%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'" + '17.71'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'" + '  -0.04%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'" + '1.008'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'" + '  -0.49%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'" + '5.926'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'" + '  -1.02%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'" + '28.555.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'" + '  +0.46%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'" + '  -0.78%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'" + '2.233'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'" + '  -1.35%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'" + '2.090.07'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'" + '  -0.77%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'" + '161.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'" + '  +0.65%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'" + '  -0.19%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'" + '2.354'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'" + '  -3.31%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'" + '126.02'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'" + '  -0.03%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'" + '0.1044'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'" + '  -2.04%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'" + '1.039'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'" + '  -0.44%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'" + '5.790'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'" + '  -1.71%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'" + '3.612'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'" + '  -0.61%  '
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'" + '0.02452'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'" + '  +0.27%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'" + '0.06549'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'" + '  +0.67%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'" + '0.2160'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'" + '  -1.11%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'" + '8.855'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'" + '  -3.83%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'" + '5.054'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'" + '  +1.60%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'" + '1.252'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'" + '  -0.71%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'" + '1.189'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'" + '  +0.22%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'" + '0.6416'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'" + '  -0.86%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'" + '  -0.49%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'" + '1.008'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'" + '  -0.48%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'" + '0.6032'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'" + '  -0.33%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'" + '12.95'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'" + '  -1.27%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'" + '3.692'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'" + '  -0.13%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'" + '2.010'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'" + '  -0.67%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'" + '1.215'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'" + '  -0.20%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'" + '122.01'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'" + '  +0.01%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'" + '1.179'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'" + '  -8.69%  '
$ws.Range("E51").Style = "Normal"
